$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 79, shifting the existing
# rows 79-138 down to 81-140 (matches the dimension change A1:R138 -> A1:R140).
$ws.Rows("79:80").Insert()

# New row 79: Vega Central Mapocho de Santiago, Camote, Primera, 2022-12-22 entry
$ws.Cells.Item(79,1).Value = 9
$ws.Cells.Item(79,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(79,3).Value = "Metropolitana"
$ws.Cells.Item(79,4).Value = 44977
$ws.Cells.Item(79,5).Value = 13
$ws.Cells.Item(79,6).Value = 100114002
$ws.Cells.Item(79,7).Value = "Camote"
$ws.Cells.Item(79,8).Value = "Sin especificar"
$ws.Cells.Item(79,9).Value = "Primera"
$ws.Cells.Item(79,10).Value = 970
$ws.Cells.Item(79,11).Value = 18000
$ws.Cells.Item(79,12).Value = 19000
$ws.Cells.Item(79,13).Value = 18500
$ws.Cells.Item(79,14).Value = "`$/caja 18 kilos"
$ws.Cells.Item(79,15).Value = "Perú"
$ws.Cells.Item(79,16).Value = 1028
$ws.Cells.Item(79,17).Value = 18
$ws.Cells.Item(79,18).Value = "Hortaliza"

# New row 80: Vega Central Mapocho de Santiago, Camote, Primera, same date
$ws.Cells.Item(80,1).Value = 9
$ws.Cells.Item(80,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(80,3).Value = "Metropolitana"
$ws.Cells.Item(80,4).Value = 44977
$ws.Cells.Item(80,5).Value = 13
$ws.Cells.Item(80,6).Value = 100114002
$ws.Cells.Item(80,7).Value = "Camote"
$ws.Cells.Item(80,8).Value = "Sin especificar"
$ws.Cells.Item(80,9).Value = "Primera"
$ws.Cells.Item(80,10).Value = 700
$ws.Cells.Item(80,11).Value = 14000
$ws.Cells.Item(80,12).Value = 15000
$ws.Cells.Item(80,13).Value = 14500
$ws.Cells.Item(80,14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(80,15).Value = "Perú"
$ws.Cells.Item(80,16).Value = 806
$ws.Cells.Item(80,17).Value = 18
$ws.Cells.Item(80,18).Value = "Hortaliza"
